$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 3193.3125
$ws.Range("I82").Value = 773.75
$ws.Range("J82").Value = 3999.8333
$ws.Range("K82").Value = 2321.25
$ws.Range("L82").Value = 11999.4999
$ws.Range("M82").Value = -1915.25
$ws.Range("N82").Value = -12811.4999

$ws.Range("H85").Value = 3193.3125
$ws.Range("I85").Value = 773.75
$ws.Range("J85").Value = 3999.8333
$ws.Range("K85").Value = 2321.25
$ws.Range("L85").Value = 11999.4999
$ws.Range("M85").Value = -917.25
$ws.Range("N85").Value = -14807.4999

$ws.Range("H112").Value = 1099.7142
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()

$ws.Range("H128").Value = 46623.6
$ws.Range("J128").Value = 46623.6
$ws.Range("L128").Value = 46623.6
$ws.Range("N128").Value = -56583.6

$ws.Range("H129").Value = 853.9091
$ws.Range("I129").Value = 666.3333
$ws.Range("J129").Value = 872.6667
$ws.Range("K129").Value = 1998.9999
$ws.Range("L129").Value = 2618.0001
$ws.Range("M129").Value = 3001.0001
$ws.Range("N129").Value = -12618.0001

$ws.Range("H137").Value = 1405.85
$ws.Range("I137").Value = 1142.2667
$ws.Range("J137").Value = 2196.6
$ws.Range("K137").Value = 3426.800099999999
$ws.Range("L137").Value = 6589.799999999999
$ws.Range("M137").Value = -876.8000999999995
$ws.Range("N137").Value = -11689.8

$ws.Range("H138").Value = 2469.0557
$ws.Range("I138").Value = 1557.8148
$ws.Range("J138").Value = 3015.8
$ws.Range("K138").Value = 4673.4444
$ws.Range("L138").Value = 9047.400000000001
$ws.Range("M138").Value = 466.5555999999997
$ws.Range("N138").Value = -19327.4

$ws.Range("H141").Value = 2050.7778
$ws.Range("I141").Value = 1507.6
$ws.Range("J141").Value = 4766.6665
$ws.Range("K141").Value = 4522.799999999999
$ws.Range("L141").Value = 14299.9995
$ws.Range("M141").Value = 657.2000000000007
$ws.Range("N141").Value = -24659.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27077.908
$ws.Range("I32").Value = 4807.4917
$ws.Range("K32").Value = 4807.4917
$ws.Range("M32").Value = -4520.4917

$ws.Range("H61").Value = 2222.75
$ws.Range("I61").Value = 1912
$ws.Range("J61").Value = 2326.3333
$ws.Range("K61").Value = 1912
$ws.Range("L61").Value = 2326.3333
$ws.Range("M61").Value = -1700
$ws.Range("N61").Value = -2750.3333

$ws.Range("H74").Value = 2614.7693
$ws.Range("I74").Value = 1505.5883
$ws.Range("J74").Value = 4709.8887
$ws.Range("K74").Value = 1505.5883
$ws.Range("L74").Value = 4709.8887
$ws.Range("M74").Value = -631.5882999999999
$ws.Range("N74").Value = -6457.8887

$ws.Range("H77").Value = 2614.7693
$ws.Range("I77").Value = 1505.5883
$ws.Range("J77").Value = 4709.8887
$ws.Range("K77").Value = 7527.941499999999
$ws.Range("L77").Value = 23549.4435
$ws.Range("M77").Value = -3159.941499999999
$ws.Range("N77").Value = -32285.4435

$ws.Range("H132").Value = 4437.5264
$ws.Range("I132").Value = 4994.385
$ws.Range("K132").Value = 14983.155
$ws.Range("M132").Value = -12453.155

$ws.Range("H136").Value = 2222.75
$ws.Range("I136").Value = 1912
$ws.Range("J136").Value = 2326.3333
$ws.Range("K136").Value = 5736
$ws.Range("L136").Value = 6978.999899999999
$ws.Range("M136").Value = -3186
$ws.Range("N136").Value = -12078.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2585.348
$ws.Range("I134").Value = 2585.348
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7756.044
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5221.044
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25651.12
$ws.Range("I31").Value = 29787.086
$ws.Range("J31").Value = 4971.2856
$ws.Range("K31").Value = 29787.086
$ws.Range("L31").Value = 4971.2856
$ws.Range("M31").Value = -29492.086
$ws.Range("N31").Value = -5561.2856

$ws.Range("H34").Value = 25651.12
$ws.Range("I34").Value = 29787.086
$ws.Range("J34").Value = 4971.2856
$ws.Range("K34").Value = 29787.086
$ws.Range("L34").Value = 4971.2856
$ws.Range("M34").Value = -29585.086
$ws.Range("N34").Value = -5375.2856

$ws.Range("H58").Value = 10250.429
$ws.Range("I58").Value = 1492.3334
$ws.Range("J58").Value = 26015
$ws.Range("K58").Value = 1492.3334
$ws.Range("L58").Value = 26015
$ws.Range("M58").Value = -1289.3334
$ws.Range("N58").Value = -26421

$ws.Range("H132").Value = 3881.077
$ws.Range("I132").Value = 4072.0557
$ws.Range("K132").Value = 12216.1671
$ws.Range("M132").Value = -9686.167099999999

$ws.Range("H136").Value = 10250.429
$ws.Range("I136").Value = 1492.3334
$ws.Range("J136").Value = 26015
$ws.Range("K136").Value = 4477.0002
$ws.Range("L136").Value = 78045
$ws.Range("M136").Value = -1927.0002
$ws.Range("N136").Value = -83145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 117.25
$ws.Range("I47").Value = 117.25
$ws.Range("K47").Value = 351.75
$ws.Range("M47").Value = 79.25

$ws.Range("H137").Value = 2695.3684
$ws.Range("I137").Value = 3666.7856
$ws.Range("J137").Value = 2128.7083
$ws.Range("K137").Value = 11000.3568
$ws.Range("L137").Value = 6386.124899999999
$ws.Range("M137").Value = -5900.356800000001
$ws.Range("N137").Value = -16586.1249

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 13968.429
$ws.Range("J52").Value = 13968.429
$ws.Range("L52").Value = 13968.429
$ws.Range("N52").Value = -14486.429

$ws.Range("H102").Value = 2889.3914
$ws.Range("I102").Value = 2909.8333
$ws.Range("J102").Value = 2815.8
$ws.Range("K102").Value = 2909.8333
$ws.Range("L102").Value = 2815.8
$ws.Range("M102").Value = -1287.8333
$ws.Range("N102").Value = -6059.8

$ws.Range("H132").Value = 2743.3572
$ws.Range("I132").Value = 2194.2632
$ws.Range("J132").Value = 3902.5557
$ws.Range("K132").Value = 6582.7896
$ws.Range("L132").Value = 11707.6671
$ws.Range("M132").Value = -4052.7896
$ws.Range("N132").Value = -16767.6671

$ws.Range("H136").Value = 13559
$ws.Range("J136").Value = 13559
$ws.Range("L136").Value = 40677
$ws.Range("N136").Value = -45777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H81").Value = 8889.333000000001
$ws.Range("J81").Value = 8889.333000000001
$ws.Range("L81").Value = 8889.333000000001
$ws.Range("N81").Value = -10885.333

$ws.Range("H82").Value = 1602.3636
$ws.Range("I82").Value = 2174.75
$ws.Range("J82").Value = 1275.2858
$ws.Range("K82").Value = 2174.75
$ws.Range("L82").Value = 1275.2858
$ws.Range("M82").Value = -1813.75
$ws.Range("N82").Value = -1997.2858

$ws.Range("H84").Value = 8889.333000000001
$ws.Range("J84").Value = 8889.333000000001
$ws.Range("L84").Value = 26667.999
$ws.Range("N84").Value = -36651.999

$ws.Range("H85").Value = 1602.3636
$ws.Range("I85").Value = 2174.75
$ws.Range("J85").Value = 1275.2858
$ws.Range("K85").Value = 2174.75
$ws.Range("L85").Value = 1275.2858
$ws.Range("M85").Value = -926.75
$ws.Range("N85").Value = -3771.2858

$ws.Range("H94").Value = 27520
$ws.Range("J94").Value = 27520
$ws.Range("L94").Value = 27520
$ws.Range("N94").Value = -28872

$ws.Range("H119").Value = 39985
$ws.Range("J119").Value = 39985
$ws.Range("L119").Value = 39985
$ws.Range("N119").Value = -49661

$ws.Range("H132").Value = 4964.1875
$ws.Range("I132").Value = 5638.5454
$ws.Range("J132").Value = 3480.6
$ws.Range("K132").Value = 16915.6362
$ws.Range("L132").Value = 10441.8
$ws.Range("M132").Value = -14385.6362
$ws.Range("N132").Value = -15501.8

$ws.Range("H136").Value = 1477.0952
$ws.Range("I136").Value = 1407.7333
$ws.Range("J136").Value = 1650.5
$ws.Range("K136").Value = 4223.199900000001
$ws.Range("L136").Value = 4951.5
$ws.Range("M136").Value = -1673.199900000001
$ws.Range("N136").Value = -10051.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 22276
$ws.Range("I55").Value = 499
$ws.Range("K55").Value = 499
$ws.Range("M55").Value = -222

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H122").Value = 2474
$ws.Range("J122").Value = 3330
$ws.Range("L122").Value = 9990
$ws.Range("N122").Value = -14890

$ws.Range("H126").Value = 1707.1428
$ws.Range("I126").Value = 1731.8636
$ws.Range("J126").Value = 1616.5
$ws.Range("K126").Value = 5195.5908
$ws.Range("L126").Value = 4849.5
$ws.Range("M126").Value = -2725.5908
$ws.Range("N126").Value = -9789.5

$ws.Range("H132").Value = 4887.625
$ws.Range("I132").Value = 6088.778
$ws.Range("J132").Value = 3343.2856
$ws.Range("K132").Value = 18266.334
$ws.Range("L132").Value = 10029.8568
$ws.Range("M132").Value = -15736.334
$ws.Range("N132").Value = -15089.8568

$ws.Range("H136").Value = 2030.909
$ws.Range("I136").Value = 956
$ws.Range("J136").Value = 2347.0588
$ws.Range("K136").Value = 2868
$ws.Range("L136").Value = 7041.176399999999
$ws.Range("M136").Value = -318
$ws.Range("N136").Value = -12141.1764
